$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1:I1 ---
$ws.Cells.Item(1,6).Value = "train_auc_mean"
$ws.Cells.Item(1,7).Value = "train_auc_std"
$ws.Cells.Item(1,8).Value = "test_auc_mean"
$ws.Cells.Item(1,9).Value = "test_auc_std"

# Copy the header style (bold, bordered, centered) from an existing header cell
$ws.Range("B1").Copy()
$ws.Range("F1:I1").PasteSpecial(-4122)

# --- Update existing B:E data values ---
$ws.Cells.Item(2,2).Value = 0.6657288237054775
$ws.Cells.Item(2,3).Value = 0.02030859421041251
$ws.Cells.Item(2,4).Value = 0.6294258373205742
$ws.Cells.Item(2,5).Value = 0.01489030269817745

$ws.Cells.Item(3,2).Value = 0.6484286141873691
$ws.Cells.Item(3,3).Value = 0.002486989293964791
$ws.Cells.Item(3,4).Value = 0.6296650717703349
$ws.Cells.Item(3,5).Value = 0.01682483346652369

$ws.Cells.Item(4,2).Value = 0.7585154145465429
$ws.Cells.Item(4,3).Value = 0.02194213185177669
$ws.Cells.Item(4,4).Value = 0.6504784688995214
$ws.Cells.Item(4,5).Value = 0.01007628693236968

$ws.Cells.Item(5,2).Value = 0.7221191260101765
$ws.Cells.Item(5,3).Value = 0.01081604556517469
$ws.Cells.Item(5,4).Value = 0.6447368421052632
$ws.Cells.Item(5,5).Value = 0.0148441072320465

$ws.Cells.Item(6,2).Value = 0.7060161628255013
$ws.Cells.Item(6,3).Value = 0.02514317068757808
$ws.Cells.Item(6,4).Value = 0.6464114832535885
$ws.Cells.Item(6,5).Value = 0.008777003027676644

# --- New F:I data values ---
$ws.Cells.Item(2,6).Value = 0.8633649170189001
$ws.Cells.Item(2,7).Value = 0.02419685626886079
$ws.Cells.Item(2,8).Value = 0.8226899596653634
$ws.Cells.Item(2,9).Value = 0.01229899196398712

$ws.Cells.Item(3,6).Value = 0.8377362221948003
$ws.Cells.Item(3,7).Value = 0.00515588624869271
$ws.Cells.Item(3,8).Value = 0.8250505286560126
$ws.Cells.Item(3,9).Value = 0.01405039152506807

$ws.Cells.Item(4,6).Value = 0.9429673530162714
$ws.Cells.Item(4,7).Value = 0.0103405237850859
$ws.Cells.Item(4,8).Value = 0.8567615011325733
$ws.Cells.Item(4,9).Value = 0.006371552262075229

$ws.Cells.Item(5,6).Value = 0.9081177555601091
$ws.Cells.Item(5,7).Value = 0.00872448205911735
$ws.Cells.Item(5,8).Value = 0.8477274261429854
$ws.Cells.Item(5,9).Value = 0.009114870178775806

$ws.Cells.Item(6,6).Value = 0.9017459306803971
$ws.Cells.Item(6,7).Value = 0.01583270590270821
$ws.Cells.Item(6,8).Value = 0.8548609096269795
$ws.Cells.Item(6,9).Value = 0.009381170360451847
